$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 75770
$ws.Range("J57").Value = 75770
$ws.Range("L57").Value = 227310
$ws.Range("N57").Value = -228308
$ws.Range("H74").Value = 6481.75
$ws.Range("I74").Value = 6713.8237
$ws.Range("K74").Value = 6713.8237
$ws.Range("M74").Value = -5777.8237
$ws.Range("H77").Value = 6481.75
$ws.Range("I77").Value = 6713.8237
$ws.Range("K77").Value = 33569.1185
$ws.Range("M77").Value = -28889.1185
$ws.Range("H98").Value = 702.21875
$ws.Range("I98").Value = 695.1
$ws.Range("K98").Value = 695.1
$ws.Range("M98").Value = 802.9
$ws.Range("H122").Value = 702.21875
$ws.Range("I122").Value = 695.1
$ws.Range("K122").Value = 2085.3
$ws.Range("M122").Value = 364.6999999999998
$ws.Range("H132").Value = 2156.875
$ws.Range("I132").Value = 2185.5386
$ws.Range("K132").Value = 6556.6158
$ws.Range("M132").Value = -4026.6158
$ws.Range("H137").Value = 10028
$ws.Range("I137").Value = 2571.6924
$ws.Range("K137").Value = 7715.0772
$ws.Range("M137").Value = -5165.0772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 480.25
$ws.Range("I4").Value = 590
$ws.Range("J4").Value = 151
$ws.Range("K4").Value = 590
$ws.Range("L4").Value = 151
$ws.Range("M4").Value = -474
$ws.Range("N4").Value = -383
$ws.Range("H61").Value = 2682.1428
$ws.Range("I61").Value = 2682.1428
$ws.Range("K61").Value = 2682.1428
$ws.Range("M61").Value = -2470.1428
$ws.Range("H132").Value = 3280
$ws.Range("I132").Value = 3334.3572
$ws.Range("J132").Value = 3110.889
$ws.Range("K132").Value = 10003.0716
$ws.Range("L132").Value = 9332.667000000001
$ws.Range("M132").Value = -7473.071599999999
$ws.Range("N132").Value = -14392.667
$ws.Range("H136").Value = 2682.1428
$ws.Range("I136").Value = 2682.1428
$ws.Range("K136").Value = 8046.428400000001
$ws.Range("M136").Value = -5496.428400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4042
$ws.Range("I134").Value = 3815.125
$ws.Range("K134").Value = 11445.375
$ws.Range("M134").Value = -8910.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2366.543
$ws.Range("I31").Value = 1669.9333
$ws.Range("J31").Value = 2889
$ws.Range("K31").Value = 1669.9333
$ws.Range("L31").Value = 2889
$ws.Range("M31").Value = -1374.9333
$ws.Range("N31").Value = -3479
$ws.Range("H34").Value = 2366.543
$ws.Range("I34").Value = 1669.9333
$ws.Range("J34").Value = 2889
$ws.Range("K34").Value = 1669.9333
$ws.Range("L34").Value = 2889
$ws.Range("M34").Value = -1467.9333
$ws.Range("N34").Value = -3293
$ws.Range("H50").Value = 32000
$ws.Range("J50").Value = 32000
$ws.Range("L50").Value = 32000
$ws.Range("N50").Value = -33250
$ws.Range("H58").Value = 2020
$ws.Range("I58").Value = 1737.8
$ws.Range("K58").Value = 1737.8
$ws.Range("M58").Value = -1534.8
$ws.Range("H132").Value = 4548.5835
$ws.Range("I132").Value = 4158.9375
$ws.Range("K132").Value = 12476.8125
$ws.Range("M132").Value = -9946.8125
$ws.Range("H134").Value = 3579
$ws.Range("I134").Value = 3696.4
$ws.Range("J134").Value = 2796.3333
$ws.Range("K134").Value = 11089.2
$ws.Range("L134").Value = 8388.999899999999
$ws.Range("M134").Value = -8554.200000000001
$ws.Range("N134").Value = -13458.9999
$ws.Range("H136").Value = 2020
$ws.Range("I136").Value = 1737.8
$ws.Range("K136").Value = 5213.4
$ws.Range("M136").Value = -2663.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 184
$ws.Range("I2").Value = 139.25
$ws.Range("K2").Value = 835.5
$ws.Range("M2").Value = -722.5
$ws.Range("H5").Value = 2175.6316
$ws.Range("J5").Value = 2225.7646
$ws.Range("L5").Value = 6677.293799999999
$ws.Range("N5").Value = -6901.293799999999
$ws.Range("H58").Value = 3500
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H116").Value = 792.25
$ws.Range("I116").Value = 792.25
$ws.Range("K116").Value = 2376.75
$ws.Range("M116").Value = 1065.25
$ws.Range("H135").Value = 2175.6316
$ws.Range("J135").Value = 2225.7646
$ws.Range("L135").Value = 20031.8814
$ws.Range("N135").Value = -25101.8814

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 671.7059
$ws.Range("I107").Value = 598.7
$ws.Range("J107").Value = 776
$ws.Range("K107").Value = 598.7
$ws.Range("L107").Value = 776
$ws.Range("M107").Value = 1321.3
$ws.Range("N107").Value = -4616
$ws.Range("H113").Value = 288285.56
$ws.Range("I113").Value = 335833.16
$ws.Range("K113").Value = 335833.16
$ws.Range("M113").Value = -333663.16
$ws.Range("H126").Value = 5093.6
$ws.Range("I126").Value = 4812
$ws.Range("J126").Value = 5124.8887
$ws.Range("K126").Value = 14436
$ws.Range("L126").Value = 15374.6661
$ws.Range("M126").Value = -11966
$ws.Range("N126").Value = -20314.6661
$ws.Range("H132").Value = 3139.3416
$ws.Range("I132").Value = 2697.0312
$ws.Range("J132").Value = 4712
$ws.Range("K132").Value = 8091.0936
$ws.Range("L132").Value = 14136
$ws.Range("M132").Value = -5561.0936
$ws.Range("N132").Value = -19196
$ws.Range("H134").Value = 35400
$ws.Range("J134").Value = 35400
$ws.Range("L134").Value = 106200
$ws.Range("N134").Value = -111270
$ws.Range("H136").Value = 40766
$ws.Range("J136").Value = 40766
$ws.Range("L136").Value = 122298
$ws.Range("N136").Value = -127398

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 19989.666
$ws.Range("J2").Value = 19989.666
$ws.Range("L2").Value = 19989.666
$ws.Range("N2").Value = -20213.666
$ws.Range("H22").Value = 4527.2856
$ws.Range("I22").Value = 4698
$ws.Range("J22").Value = 4299.6665
$ws.Range("K22").Value = 4698
$ws.Range("L22").Value = 4299.6665
$ws.Range("M22").Value = -4403
$ws.Range("N22").Value = -4889.6665
$ws.Range("H27").Value = 4527.2856
$ws.Range("I27").Value = 4698
$ws.Range("J27").Value = 4299.6665
$ws.Range("K27").Value = 4698
$ws.Range("L27").Value = 4299.6665
$ws.Range("M27").Value = -4591
$ws.Range("N27").Value = -4513.6665
$ws.Range("H55").Value = 348.58334
$ws.Range("I55").Value = 849.5
$ws.Range("K55").Value = 849.5
$ws.Range("M55").Value = -676.5
$ws.Range("H132").Value = 7942199.5
$ws.Range("I132").Value = 10105799
$ws.Range("J132").Value = 8999.666999999999
$ws.Range("K132").Value = 30317397
$ws.Range("L132").Value = 26999.001
$ws.Range("M132").Value = -30314867
$ws.Range("N132").Value = -32059.001
$ws.Range("H136").Value = 5292412.5
$ws.Range("I136").Value = 5849393
$ws.Range("K136").Value = 17548179
$ws.Range("M136").Value = -17545629

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 40330
$ws.Range("J94").Value = 40330
$ws.Range("L94").Value = 40330
$ws.Range("N94").Value = -42132
$ws.Range("H132").Value = 5684.886
$ws.Range("I132").Value = 4956.628
$ws.Range("K132").Value = 14869.884
$ws.Range("M132").Value = -12339.884
$ws.Range("H136").Value = 4001519.5
$ws.Range("I136").Value = 5001435
$ws.Range("K136").Value = 15004305
$ws.Range("M136").Value = -15001755
